$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: remove trailing maqaf from A17
$ws.Range("A17").Value = "מִיּֽוֹם"

# Rewrite verse 11 (rows 94-99) with new Hebrew/Russian text and new styling
$ws.Range("A94").Value = "יִרְעַם"
$ws.Range("B94").Value = "Пусть гремит"
$ws.Range("A95").Value = "הַיָּם"
$ws.Range("B95").Value = "Море"
$ws.Range("A96").Value = "וּמְלֹאוֹ"
$ws.Range("B96").Value = "И все"
$ws.Range("A97").Value = "תֵּבֵל"
$ws.Range("B97").Value = "Вселенная"
$ws.Range("A98").Value = "וְיֹשְׁבֵי"
$ws.Range("B98").Value = "И живущие"
$ws.Range("A99").Value = "בָהּ"
$ws.Range("B99").Value = "В ней"

# Apply new font/alignment styling to A94:A99 (right aligned) and B94:B99 (left aligned)
$aRange = $ws.Range("A94:A99")
$aRange.Font.Name = "Calibri"
$aRange.Font.Size = 10
$aRange.Font.Color = 0
$aRange.HorizontalAlignment = -4152

$bRange = $ws.Range("B94:B99")
$bRange.Font.Name = "Calibri"
$bRange.Font.Size = 10
$bRange.Font.Color = 0
$bRange.HorizontalAlignment = -4131

# Remove the now-extra row (old row 100 held the last part of the old verse 11 text)
$ws.Rows.Item(100).Delete()
